# "Added more visual ideas"
# - Rename the sheet from "Sheet2" to "Draft"
# - Move the active selection on the sheet from C22 to B14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Draft"

# Update the current selection/active cell on the sheet
$ws.Range("B14").Select()
